$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New test case row (row 39), mirroring the formatting of the row above it
# (border + default alignment on A/B/D/E, wrap-text on the description column C).
$ws.Range("A37:E37").Copy()
$ws.Range("A39:E39").PasteSpecial(-4122)

# Populate the new row's values (order matches how the shared strings were
# appended: description, test case id, then jira id).
$ws.Range("C39").Value = "Verify that a user's public watchlist is not visible to another user once that particular watchlist is deleted."
$ws.Range("A39").Value = "TestCase_E38"
$ws.Range("B39").Value = "OPQA-1105"
$ws.Range("D39").Value = "Y"
$ws.Range("E39").Value = "PASS"

# Move the active selection to the newly added row, matching the saved view state.
$ws.Range("D38").Select()
